# retrait de express et ajout de nextjs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column C (Sujet) / D (Exercice) content for rows 2-11: swap Express/Node
# content for TypeScript/React links rooted under new "typescript/" and
# "react/" folders. Rows 12-31 lose their C/D content entirely (those topics
# - Examen/Projet/Retour formatif rows - no longer carry a markdown link).

$ws.Range("C2").Value = "[Introduction à TypeScript](typescript/introduction_typescript.md)"
$ws.Range("D2").ClearContents()

$ws.Range("C3").Value = "[TypeScript 2](typescript/typescript_2.md)"
$ws.Range("D3").ClearContents()

$ws.Range("C4").Value = "[JavaScript asynchrone](typescript/javascript_async.md)"
$ws.Range("D4").ClearContents()

$ws.Range("C5").Value = "[Introduction à React](react/introduction_react.md)"
$ws.Range("D5").Value = "[Exercice 9 - React](react/exercice9_react.md)"

$ws.Range("C6").Value = "[React et les styles](react/react_styles.md)"
$ws.Range("D6").Value = "[Exercice 10 - React et styles](react/exercice10_react_styles.md)"

$ws.Range("C7").Value = "[Routes, contexte et API](react/react3.md)"
$ws.Range("D7").Value = "[Exercice 11 - Contexte et API](react/exercice11_context.md)"

$ws.Range("C8").Value = "[Internationalisation](react/internationalisation.md) <br />[Accessibilité](react/accessibilite.md)"
$ws.Range("D8").Value = "[Exercice 12 - Internatialisation](react/exercice12_internationalisation.md)"

$ws.Range("C9").Value = "[Authentification](react/authentification.md)"
$ws.Range("D9").Value = "[Exercice 13 - Authentification](react/exercice13_authentification.md)"

$ws.Range("C10").Value = "[React et PWA](react/pwa.md)"
$ws.Range("D10").Value = "[Exercice 14 - PWA](react/exercice14_pwa.md)"

$ws.Range("C11").Value = "[Netlify](react/netlify.md)"
$ws.Range("D11").Value = "[Exercice 15 - Netlify](react/exercice15_netlify.md)"

# Rows 12-31: remove any leftover C/D content (evaluations/projects no longer
# reference the old Express markdown pages).
for ($r = 12; $r -le 31; $r++) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("D$r").ClearContents()
}

# --- Column C width: widen from bestFit 25.66 to a fixed 56, no bestFit.
# (55.17 is the ColumnWidth value that round-trips to the stored XML width of
# exactly 56 given this font's max digit width.)
$ws.Columns.Item(3).ColumnWidth = 55.17

# --- Selection moves from C20 to C5.
$ws.Range("C5").Select()
